$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

$ws.Range("D2").Value = "26.679.77"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.642.79"
$ws.Range("E3").Value = "  +0.89%  "
Set-TextValue $ws "D4" "1.01"
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue $ws "D5" "214.97"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  +1.56%  "
Set-TextValue $ws "D7" "1.01"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("E9").Value = "  +0.84%  "
Set-TextValue $ws "D10" "19.06"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.871.42"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws "D13" "4.18"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.627.34"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  +1.47%  "
Set-TextValue $ws "D16" "65.08"
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "26.689.90"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "0.0₃0743"
$ws.Range("E18").Value = "  +0.51%  "
Set-TextValue $ws "D19" "216.34"
$ws.Range("E19").Value = "  +0.41%  "
Set-TextValue $ws "D20" "1.01"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +0.90%  "
Set-TextValue $ws "D22" "6.26"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +1.94%  "
Set-TextValue $ws "D24" "2.25"
$ws.Range("E24").Value = "  +13.74%  "
Set-TextValue $ws "D25" "145.44"
$ws.Range("E25").Value = "  -2.13%  "
Set-TextValue $ws "D26" "1.01"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -0.09%  "
Set-TextValue $ws "D28" "7.16"
$ws.Range("E28").Value = "  +4.52%  "
Set-TextValue $ws "D29" "15.70"
$ws.Range("E29").Value = "  +0.96%  "
Set-TextValue $ws "D30" "0.0515"
$ws.Range("E30").Value = "  +1.47%  "
$ws.Range("E31").Value = "  +1.32%  "
Set-TextValue $ws "D32" "3.36"
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("E33").Value = "  +2.98%  "
$ws.Range("D34").Value = "1.278.50"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("E36").Value = "  +1.26%  "
Set-TextValue $ws "D37" "0.0178"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("E38").Value = "  +6.45%  "
Set-TextValue $ws "D39" "0.828"
$ws.Range("E39").Value = "  +3.81%  "
Set-TextValue $ws "D40" "1.01"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("D44").Value = "1.781.67"
$ws.Range("E44").Value = "  +0.92%  "
Set-TextValue $ws "D45" "91.76"
$ws.Range("E45").Value = "  -0.84%  "
Set-TextValue $ws "D46" "59.18"
$ws.Range("E46").Value = "  +7.91%  "
Set-TextValue $ws "D47" "1.59"
$ws.Range("E47").Value = "  +1.46%  "
$ws.Range("E48").Value = "  +0.95%  "
Set-TextValue $ws "D49" "7.74"
$ws.Range("E49").Value = "  +1.10%  "
Set-TextValue $ws "D50" "0.0964"
$ws.Range("E50").Value = "  +1.83%  "
Set-TextValue $ws "D51" "0.407"
$ws.Range("E51").Value = "  -0.59%  "
